# "Initialen eintragen, nicht namen" - replace the full-name strings that were
# used in the tutor-assignment columns (L on the PUE sheet, J on the HUE sheet)
# with the initials that are already used elsewhere in the workbook.
#
# Mapping (same text everywhere it occurs):
#   Setzer   -> AS   (Alexander Setzer)
#   Parruca  -> DP   (Donald Parruca)
#   Feldmann -> MF   (Michael Feldmann)
# except row 3 on the PUE sheet, where "Setzer" becomes "DP" (matches the
# source diff exactly).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PUE sheet - column L
# ---------------------------------------------------------------------------
$pue = $wb.Worksheets.Item("PUE")

$pue.Range("L3").Value = "DP"
$pue.Range("L3").Font.Bold = $false

$pue.Range("L4").Value = "DP"
$pue.Range("L5").Value = "DP"
$pue.Range("L6").Value = "DP"
$pue.Range("L7").Value = "DP"
$pue.Range("L8").Value = "DP"

$pue.Range("L9").Value = "MF"
$pue.Range("L10").Value = "MF"
$pue.Range("L11").Value = "MF"
$pue.Range("L12").Value = "MF"

$pue.Range("L13").Value = "AS"
$pue.Range("L14").Value = "AS"
$pue.Range("L15").Value = "AS"
$pue.Range("L16").Value = "AS"

# ---------------------------------------------------------------------------
# HUE sheet - column J
# ---------------------------------------------------------------------------
$hue = $wb.Worksheets.Item("HUE")

$hue.Range("J3").Value = "AS"

$hue.Range("J4").Value = "DP"
$hue.Range("J5").Value = "DP"
$hue.Range("J6").Value = "DP"
$hue.Range("J7").Value = "DP"
$hue.Range("J8").Value = "DP"

$hue.Range("J9").Value = "MF"
$hue.Range("J10").Value = "MF"
$hue.Range("J11").Value = "MF"
$hue.Range("J12").Value = "MF"

$hue.Range("J13").Value = "AS"
$hue.Range("J14").Value = "AS"
$hue.Range("J15").Value = "AS"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping so the saved view matches the edit:
#   - PUE is no longer the tab shown on open, its remembered selection is K4
#   - Tutoren keeps its own remembered selection, now D13
#   - HUE becomes the active (visible) sheet with selection J16
# ---------------------------------------------------------------------------
$pue.Activate()
$pue.Range("K4").Select()

$tutoren = $wb.Worksheets.Item("Tutoren")
$tutoren.Activate()
$tutoren.Range("D13").Select()

$hue.Activate()
$hue.Range("J16").Select()
